$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A8 with combined tuple-like strings, consolidating rows 2-25 into 2-8
$ws.Range("A2").Value = "('Disenchant', ['{1}{W}', 'Instant', 'Destroy target artifact or enchantment.'])"
$ws.Range("A3").Value = "('Fireball', ['{X}{R}', 'Sorcery', 'This spell costs {1} more to cast for each target beyond the first.', 'Fireball deals X damage divided evenly, rounded down, among any number of targets.'])"
$ws.Range("A4").Value = "('Forest', ['Basic Land — Forest', '({T}: Add {G}.)'])"
$ws.Range("A5").Value = "('Island', ['Basic Land — Island', '({T}: Add {U}.)'])"
$ws.Range("A6").Value = "('Mountain', ['Basic Land — Mountain', '({T}: Add {R}.)'])"
$ws.Range("A7").Value = "('Plains', ['Basic Land — Plains', '({T}: Add {W}.)'])"
$ws.Range("A8").Value = "('Swamp', ['Basic Land — Swamp', '({T}: Add {B}.)'])"

# Remove now-unused rows 9-25 so the sheet dimension shrinks to A1:A8
$ws.Range("A9:A25").EntireRow.Delete()
